$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 64
$ws1.Range("F7").Value = 2658
$ws1.Range("F8").Value = 1164
$ws1.Range("F9").Value = 247
$ws1.Range("F10").Value = 102
$ws1.Range("F11").Value = 9771
$ws1.Range("F13").Value = 244
$ws1.Range("F14").Value = 1
$ws1.Range("F15").Value = 596
$ws1.Range("F16").Value = 11678
$ws1.Range("F17").Value = 11978

# Sheet "全部类型" (sheet4): update "想去人数" (column F) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 64
$ws4.Range("F7").Value = 2658
$ws4.Range("F9").Value = 1164
$ws4.Range("F10").Value = 247
$ws4.Range("F11").Value = 102
$ws4.Range("F12").Value = 9771
$ws4.Range("F14").Value = 244
$ws4.Range("F15").Value = 1
$ws4.Range("F16").Value = 596
$ws4.Range("F17").Value = 11678
$ws4.Range("F18").Value = 11978

$wb.Save()
